$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")
$ws.Range("A2").Value = "tavalinetont48"
$ws.Range("C2").Value = "puhtaloom48"
$ws.Range("E2").Value = "filmweird48"
